$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15
$ws.Range("C2").Formula = "=1.2*B2"
$ws.Range("D2").Value = 0.82
$ws.Range("G2").Value = 0.5

$ws.Range("D2").Select()
